# Generate Report for Handback
# Marks the two handed-off files as handed back (status text + new
# "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime"
# entries on each locale sheet), and widens a few columns to fit the new
# (longer) content - mirroring the localization-status report regeneration.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$mdFile1 = "01bbae0a-7cbe-451a-9851-9b39a347bbca.md"
$mdFile2 = "212c757b-9df1-46f1-8d2b-b5f317d9d6a2.md"

$url1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/807174ebae48cd4cc1eaa2808875d0f2f428c5eb/e2e/01bbae0a-7cbe-451a-9851-9b39a347bbca.md"
$url2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/807174ebae48cd4cc1eaa2808875d0f2f428c5eb/e2e/212c757b-9df1-46f1-8d2b-b5f317d9d6a2.md"

# ---------------------------------------------------------------------
# Overview sheet: status columns (zh-cn / de-de) now read "handed back"
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# Widen the now-longer status columns to fit.
$overview.Columns.Item(5).ColumnWidth = 29.15
$overview.Columns.Item(6).ColumnWidth = 29.15

# ---------------------------------------------------------------------
# Per-locale detail sheets
# ---------------------------------------------------------------------
$locales = @(
    @{ Name = "zh-cn"; Xlf1 = "01bbae0a-7cbe-451a-9851-9b39a347bbca.1c60061316b2676ae2aca24389963216ca9cc6f7.zh-cn.xlf"; Xlf2 = "212c757b-9df1-46f1-8d2b-b5f317d9d6a2.5f44469879b8c2a7e1f0d5c85a5ab8d5476572e0.zh-cn.xlf"; HandbackTime = "2016-09-07 11:27:42" },
    @{ Name = "de-de"; Xlf1 = "01bbae0a-7cbe-451a-9851-9b39a347bbca.1c60061316b2676ae2aca24389963216ca9cc6f7.de-de.xlf"; Xlf2 = "212c757b-9df1-46f1-8d2b-b5f317d9d6a2.5f44469879b8c2a7e1f0d5c85a5ab8d5476572e0.de-de.xlf"; HandbackTime = "2016-09-07 11:27:51" }
)

foreach ($locale in $locales) {
    $ws = $wb.Worksheets.Item($locale.Name)

    # Status column now reflects the handback.
    $ws.Range("C2").Value = $newStatus
    $ws.Range("C3").Value = $newStatus

    # Newly generated handback artefacts: target (source) file, handback
    # xliff file and the handback timestamp, for each of the two rows.
    $ws.Range("I2").Value = $mdFile1
    $ws.Range("J2").Value = $locale.Xlf1
    $ws.Range("K2").Value = $locale.HandbackTime

    $ws.Range("I3").Value = $mdFile2
    $ws.Range("J3").Value = $locale.Xlf2
    $ws.Range("K3").Value = $locale.HandbackTime

    # Re-create the hyperlinks so that both the source-file link (column A)
    # and the newly-added target-file link (column I) are present, in row
    # order, for each of the two data rows.
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), $url1, "", "", $mdFile1)
    $ws.Hyperlinks.Add($ws.Range("I2"), $url1, "", "", $mdFile1)
    $ws.Hyperlinks.Add($ws.Range("A3"), $url2, "", "", $mdFile2)
    $ws.Hyperlinks.Add($ws.Range("I3"), $url2, "", "", $mdFile2)

    # Widen the status/target/handback columns to fit the new content.
    $ws.Columns.Item(3).ColumnWidth = 29.15
    $ws.Columns.Item(9).ColumnWidth = 39.15
    $ws.Columns.Item(10).ColumnWidth = 39.15
}
